$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.5
$summary.Range("C2").Value = 0.5
$summary.Range("D2").Value = 1
$summary.Range("E2").Value = 0.6666666666666666
$summary.Range("F2").Value = 0.8333333333333334
$summary.Range("G2").Value = 0.9629629629629629
$summary.Range("H2").Value = 0.7561966081723687
$summary.Range("I2").Value = 534
$summary.Range("J2").Value = 534
$summary.Range("K2").Value = 0
$summary.Range("L2").Value = 0

$classRep = $wb.Worksheets.Item("Classification Report")
$classRep.Range("B2").Value = 0
$classRep.Range("C2").Value = 0
$classRep.Range("D2").Value = 0

$classRep.Range("B3").Value = 0.5
$classRep.Range("C3").Value = 1
$classRep.Range("D3").Value = 0.6666666666666666

$classRep.Range("B4").Value = 0.5
$classRep.Range("C4").Value = 0.5
$classRep.Range("D4").Value = 0.5
$classRep.Range("E4").Value = 0.5

$classRep.Range("B5").Value = 0.25
$classRep.Range("C5").Value = 0.5
$classRep.Range("D5").Value = 0.3333333333333333

$classRep.Range("B6").Value = 0.25
$classRep.Range("C6").Value = 0.5
$classRep.Range("D6").Value = 0.3333333333333333

$confMat = $wb.Worksheets.Item("Confusion Matrix")
$confMat.Range("B2").Value = 0
$confMat.Range("C2").Value = 534

$confMat.Range("B3").Value = 0
$confMat.Range("C3").Value = 534
